$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("번역")
$ws2 = $wb.Worksheets.Item("용어통일")

# Order matches original authoring order so new shared-string entries land
# in the same sequence as the target workbook.
$ws2.Range("A10").Value = "테일러 급수, 로랑 급수"
$ws2.Range("A11").Value = "근의 분류"
$ws2.Range("A12").Value = "영점, 극점, 본질적 특이점, n차 극점, 단순 극점"

# Sheet "번역" (translation) - add column C/D values to row 75
$ws1.Range("C75").Value = "수렴 반지름"
$ws1.Range("D75").Value = "수렴 반지름으로 통일"

$ws2.Range("A13").Value = "수렴 반지름"

# Update selections to match final state
$ws1.Range("D75").Select()
$ws2.Range("A14").Select()
